$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.182.60"
$ws.Range("E2").Value = "  +3.60%  "

# Row 3
$ws.Range("D3").Value = "3.271.68"
$ws.Range("E3").Value = "  +2.88%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.68"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.84"
$ws.Range("E6").Value = "  +6.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +0.63%  "

# Row 9
$ws.Range("D9").Value = "3.263.29"
$ws.Range("E9").Value = "  +2.60%  "

# Row 10
$ws.Range("E10").Value = "  +7.43%  "

# Row 11
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("E12").Value = "  +6.31%  "

# Row 13
$ws.Range("D13").Value = "3.833.50"
$ws.Range("E13").Value = "  +3.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.137"
$ws.Range("E14").Value = "  +0.88%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.56"
$ws.Range("E15").Value = "  +4.18%  "

# Row 16
$ws.Range("D16").Value = "68.087.66"
$ws.Range("E16").Value = "  +3.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000170"
$ws.Range("E17").Value = "  +3.65%  "

# Row 18
$ws.Range("D18").Value = "3.271.55"
$ws.Range("E18").Value = "  +3.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  +2.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.54"
$ws.Range("E20").Value = "  +4.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.36"
$ws.Range("E21").Value = "  +4.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.66"
$ws.Range("E22").Value = "  +5.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.32"
$ws.Range("E24").Value = "  +2.93%  "

# Row 25
$ws.Range("E25").Value = "  +3.77%  "

# Row 26
$ws.Range("E26").Value = "  +4.74%  "

# Row 27
$ws.Range("E27").Value = "  -2.17%  "

# Row 28
$ws.Range("E28").Value = "  +2.23%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.07%  "

# Row 30
$ws.Range("E30").Value = "  +3.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.71"
$ws.Range("E31").Value = "  +5.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.87"
$ws.Range("E32").Value = "  +3.99%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  +5.85%  "

# Row 35
$ws.Range("E35").Value = "  +4.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  +5.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.93"
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +1.77%  "

# Row 39
$ws.Range("E39").Value = "  +2.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.81"
$ws.Range("E40").Value = "  +11.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.88"
$ws.Range("E41").Value = "  +1.77%  "

# Row 42
$ws.Range("E42").Value = "  +10.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  +4.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.78"
$ws.Range("E44").Value = "  +6.64%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "351.44"
$ws.Range("E45").Value = "  +6.55%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.688.25"
$ws.Range("E46").Value = "  +1.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.87"
$ws.Range("E47").Value = "  +3.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0682"
$ws.Range("E48").Value = "  +3.43%  "

# Row 49
$ws.Range("E49").Value = "  +2.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  +5.85%  "

# Row 51
$ws.Range("E51").Value = "  +0.50%  "
